$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.944.35'
$ws.Range("E2").Value = '  +0.58%  '

$ws.Range("D3").Value = '2.788.91'
$ws.Range("E3").Value = '  -1.32%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '359.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.25%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.566'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.24%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.595'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.12%  '

$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("E12").Value = '  +0.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.55'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.22%  '

$ws.Range("D15").Value = '3.223.19'
$ws.Range("E15").Value = '  -1.52%  '

$ws.Range("D16").Value = '2.792.89'
$ws.Range("E16").Value = '  -1.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.934'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.58%  '

$ws.Range("D18").Value = '51.860.60'
$ws.Range("E18").Value = '  +0.75%  '

$ws.Range("E19").Value = '  +0.92%  '

$ws.Range("E20").Value = '  -0.98%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.49%  '

$ws.Range("E22").Value = '  -1.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '274.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.37%  '

$ws.Range("E24").Value = '  +0.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.99%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.20%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.05%  '

$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.145'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.78%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0468'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.02%  '

$ws.Range("E32").Value = '  +1.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '33.97'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.54%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0843'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.46%  '

$ws.Range("E36").Value = '  +6.54%  '

$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.23'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.11'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.40%  '

$ws.Range("E40").Value = '  -3.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.58%  '

$ws.Range("E42").Value = '  -1.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '122.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.19%  '

$ws.Range("E44").Value = '  -2.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.09'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.85%  '

$ws.Range("D46").Value = '2.076.55'
$ws.Range("E46").Value = '  -0.10%  '

$ws.Range("E47").Value = '  -1.96%  '

$ws.Range("E48").Value = '  -5.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.71'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.937'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.28%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.73%  '
